$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $ws.Range("ZZ1").Value = "'" + $text
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "42.974.22"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.302.77"
$ws.Range("E3").Value = "  -0.02%  "
Set-TextValue "D5" "304.83"
$ws.Range("E5").Value = "  +1.27%  "
Set-TextValue "D6" "97.78"
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -1.85%  "
$ws.Range("E8").Value = "  +0.01%  "
Set-TextValue "D9" "0.506"
$ws.Range("E9").Value = "  -1.97%  "
Set-TextValue "D10" "35.76"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("E11").Value = "  -0.04%  "
Set-TextValue "D12" "18.15"
$ws.Range("E12").Value = "  +0.77%  "
Set-TextValue "D13" "0.118"
$ws.Range("E13").Value = "  +1.23%  "
Set-TextValue "D14" "6.78"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "2.663.06"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "2.299.51"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "42.895.79"
$ws.Range("E18").Value = "  -0.17%  "
Set-TextValue "D19" "12.63"
$ws.Range("E19").Value = "  -5.56%  "
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -0.22%  "
Set-TextValue "D21" "6.04"
$ws.Range("E21").Value = "  -1.16%  "
Set-TextValue "D22" "67.99"
$ws.Range("E22").Value = "  -0.42%  "
Set-TextValue "D23" "236.80"
$ws.Range("E23").Value = "  -0.65%  "
Set-TextValue "D25" "2.48"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -0.33%  "
Set-TextValue "D28" "25.45"
$ws.Range("E28").Value = "  +3.01%  "
Set-TextValue "D29" "167.53"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  +0.43%  "
Set-TextValue "D32" "33.18"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +0.13%  "
Set-TextValue "D34" "4.81"
$ws.Range("E34").Value = "  +0.27%  "
Set-TextValue "D35" "5.02"
$ws.Range("E35").Value = "  -2.94%  "
Set-TextValue "D36" "17.34"
$ws.Range("E36").Value = "  -4.09%  "
$ws.Range("E37").Value = "  -1.19%  "
Set-TextValue "D38" "0.0690"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  -1.32%  "
Set-TextValue "D40" "1.76"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").Value = "2.008.39"
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  -2.16%  "
Set-TextValue "D45" "10.05"
$ws.Range("E45").Value = "  -1.34%  "
Set-TextValue "D46" "17.95"
$ws.Range("E46").Value = "  +3.67%  "
Set-TextValue "D47" "2.10"
$ws.Range("E47").Value = "  -3.13%  "
Set-TextValue "D48" "2.79"
$ws.Range("E48").Value = "  -1.22%  "
Set-TextValue "D49" "2.89"
$ws.Range("E49").Value = "  +2.69%  "
Set-TextValue "D50" "53.69"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "2.530.61"
$ws.Range("E51").Value = "  +0.06%  "

$ws.Range("ZZ1").Clear()
